# Updated cryptos list on Sat Aug 24 07:35:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.054.20"
$ws.Range("E2").Value = "  +5.14%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.733.43"
$ws.Range("E3").Value = "  +2.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.74"
$ws.Range("E5").Value = "  -0.24%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.46"
$ws.Range("E6").Value = "  +6.42%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.87%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.760.47"
$ws.Range("E9").Value = "  +3.31%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +1.91%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +4.73%  "

# Row 12 - was TRON, becomes Cardano
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  +2.16%  "

# Row 13 - was Cardano, becomes TRON
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.161"
$ws.Range("E13").Value = "  +4.60%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.233.12"
$ws.Range("E14").Value = "  +2.80%  "

# Row 15 - Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.30"
$ws.Range("E15").Value = "  +2.16%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "64.059.63"
$ws.Range("E16").Value = "  +5.15%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +6.11%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.761.13"
$ws.Range("E18").Value = "  +3.24%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +3.01%  "

# Row 20 - Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.87"
$ws.Range("E20").Value = "  +2.63%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.88"
$ws.Range("E21").Value = "  +2.70%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.52%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.09%  "

# Row 24 - Polygon
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.532"
$ws.Range("E24").Value = "  -0.27%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.39"
$ws.Range("E25").Value = "  +3.63%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +4.77%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +4.64%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.05%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +12.62%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.09%  "

# Row 31 - Aptos
$ws.Range("E31").Value = "  +4.48%  "

# Row 32 - was Monero, becomes Fetch.AI
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("E32").Value = "  +16.99%  "

# Row 33 - was Fetch.AI, becomes Monero
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "171.91"
$ws.Range("E33").Value = "  +3.24%  "

# Row 34 - USDe
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.06%  "

# Row 35 - EthereumClassic
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.43"
$ws.Range("E35").Value = "  +2.69%  "

# Row 36 - NEARProtocol
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.80"
$ws.Range("E36").Value = "  +7.85%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +8.32%  "

# Row 38 - SuiNetwork
$ws.Range("E38").Value = "  +17.83%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +9.77%  "

# Row 40 - Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "343.30"
$ws.Range("E40").Value = "  +4.15%  "

# Row 41 - Filecoin
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  +5.30%  "

# Row 42 - OKB
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.40"
$ws.Range("E42").Value = "  +2.69%  "

# Row 43 - RenderToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.61"
$ws.Range("E43").Value = "  +7.25%  "

# Row 44 - EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.79"
$ws.Range("E44").Value = "  +6.03%  "

# Row 45 - was Hedera, becomes InjectiveProtocol
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.74"
$ws.Range("E45").Value = "  +6.26%  "

# Row 46 - was InjectiveProtocol, becomes Hedera
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0591"
$ws.Range("E46").Value = "  +5.25%  "

# Row 47 - Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.649"
$ws.Range("E47").Value = "  +5.38%  "

# Row 48 - Aave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.11"
$ws.Range("E48").Value = "  +3.34%  "

# Row 49 - VeChain
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0255"
$ws.Range("E49").Value = "  +3.07%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +0.72%  "

# Row 51 - FirstDigitalUSD
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -0.26%  "
